# Add a new "2022-Q3" sheet (fund holdings detail) as the 2nd sheet (right
# after "总计"), shifting "2022-Q2" ... "2020-Q4" one position to the right,
# and insert a corresponding new leading data row on the "总计" summary
# sheet (shifting its existing rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 at row 2,
#    pushing the existing 2022-Q2 .. 2020-Q4 rows down one row each.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the new bottom row (9) the same formatting as row 8's A-column cell
# (bold/centered/bordered "index" style) before shifting values into it.
$summary.Range("A8").Copy($summary.Range("A9"))

# Shift existing data rows down by one, bottom-up so we never clobber a
# row before it has been copied.
for ($r = 8; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $summary.Range("A" + $dstRow).Value = $summary.Range("A" + $srcRow).Value2
    $summary.Range("B" + $dstRow).Value = $summary.Range("B" + $srcRow).Value2
    $summary.Range("C" + $dstRow).Value = $summary.Range("C" + $srcRow).Value2
    $summary.Range("D" + $dstRow).Value = $summary.Range("D" + $srcRow).Value2
}

# Fill in the brand-new row 2 with the 2022-Q3 totals.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 0.18

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" (i.e. before
#    the worksheet currently named "2022-Q2"), then populate it with the
#    fund holdings detail rows.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($oldQ2)
$newSheet.Name = "2022-Q3"

# Pull over the header row's text + styling (bold/centered/bordered) from
# the existing quarter sheet so formatting matches the rest of the workbook.
$oldQ2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Likewise seed the "index" style (column A) used on every data row.
$oldQ2.Range("A2").Copy($newSheet.Range("A2"))

# idx, fund code, fund name, fund size, stock position, position pct, held value, rank
$rows = @(
    @(0, "562500", "华夏中证机器人ETF",             "1.57", "99.51", "3.59", "0.0564", 7),
    @(1, "159770", "天弘中证机器人ETF",             "1.04", "99.68", "3.59", "0.0373", 7),
    @(2, "167506", "安信深圳科技指数（LOF）A",       "0.82", "93.20", "3.98", "0.0326", 7),
    @(3, "562360", "银华中证机器人ETF",             "0.77", "97.23", "3.52", "0.0271", 7),
    @(4, "167507", "安信深圳科技指数（LOF）C",       "0.30", "93.20", "3.98", "0.0119", 7),
    @(5, "002068", "东方多策略灵活配置混合C",       "0.26", "55.14", "2.74", "0.0071", 3),
    @(6, "002023", "红塔红土稳健回报灵活配置混合A", "0.09", "62.51", "2.89", "0.0026", 10),
    @(7, "400023", "东方多策略灵活配置混合A",       "0.03", "55.14", "2.74", "0.0008", 3),
    @(8, "002024", "红塔红土稳健回报灵活配置混合C", "0.00", "62.51", "2.89", "0",      10)
)

$rowNum = 2
foreach ($rowData in $rows) {
    if ($rowNum -gt 2) {
        $newSheet.Range("A2").Copy($newSheet.Range("A" + $rowNum))
    }
    $newSheet.Range("A" + $rowNum).Value = $rowData[0]

    # Columns B..G carry text values (fund codes / formatted numbers kept
    # as strings, matching every other quarter sheet in this workbook).
    $textRange = $newSheet.Range("B" + $rowNum + ":G" + $rowNum)
    $textRange.NumberFormat = "@"
    $newSheet.Range("B" + $rowNum).Value = $rowData[1]
    $newSheet.Range("C" + $rowNum).Value = $rowData[2]
    $newSheet.Range("D" + $rowNum).Value = $rowData[3]
    $newSheet.Range("E" + $rowNum).Value = $rowData[4]
    $newSheet.Range("F" + $rowNum).Value = $rowData[5]
    $newSheet.Range("G" + $rowNum).Value = $rowData[6]
    $textRange.Style = "Normal"

    $newSheet.Range("H" + $rowNum).Value = $rowData[7]

    $rowNum = $rowNum + 1
}

# Last row's "held value" (G10) is a true zero, stored as a number rather
# than the text "0.0000" used elsewhere (matches how the rest of the
# workbook represents exact-zero holdings).
$newSheet.Range("G10").NumberFormat = "General"
$newSheet.Range("G10").Value = 0
$newSheet.Range("G10").Style = "Normal"

# Match the page margins used by the other quarter worksheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$summary.Select()
